$wb = $excel.ActiveWorkbook

# Rename sheets (order: GNG_TO, NB_TO, RS_TO, TOL_TO, vSAT_TO)
$wsGNG  = $wb.Worksheets.Item(1)
$wsNB   = $wb.Worksheets.Item(2)
$wsRS   = $wb.Worksheets.Item(3)
$wsTOL  = $wb.Worksheets.Item(4)
$wsvSAT = $wb.Worksheets.Item(5)

$wsGNG.Name  = "GNG_TO-16509960993353374"
$wsNB.Name   = "NB_TO-16509961008953276"
$wsRS.Name   = "RS_TO-16509961008953276"
$wsTOL.Name  = "TOL_TO-16509961009432955"
$wsvSAT.Name = "vSAT_TO-16509961010152953"

# Sheet1: GNG_TO - update B2:B5
$wsGNG.Range("B2").Value = "go_stims-1650996099295331.csv"
$wsGNG.Range("B3").Value = "GNG_stims-1650996099319296.csv"
$wsGNG.Range("B4").Value = "go_stims-1650996099319296.csv"
$wsGNG.Range("B5").Value = "GNG_stims-16509960993353374.csv"

# Sheet2: NB_TO - update B2:B10
$wsNB.Range("B2").Value = "TB-16509961008713443.csv"
$wsNB.Range("B3").Value = "ZB-match_0-16509960998232946.csv"
$wsNB.Range("B4").Value = "TB-16509961008153298.csv"
$wsNB.Range("B5").Value = "OB-16509960998632956.csv"
$wsNB.Range("B6").Value = "TB-16509961004393268.csv"
$wsNB.Range("B7").Value = "OB-16509961002233303.csv"
$wsNB.Range("B8").Value = "OB-16509961002713284.csv"
$wsNB.Range("B9").Value = "ZB-match_3-16509960995593016.csv"
$wsNB.Range("B10").Value = "ZB-match_6-1650996099359331.csv"

# Sheet3: RS_TO - no data changes

# Sheet4: TOL_TO - update B2:B7
$wsTOL.Range("B2").Value = "MM_stims-16509961009113326.csv"
$wsTOL.Range("B3").Value = "ZM_stims-16509961008953276.csv"
$wsTOL.Range("B4").Value = "MM_stims-1650996100927329.csv"
$wsTOL.Range("B5").Value = "ZM_stims-16509961009113326.csv"
$wsTOL.Range("B6").Value = "MM_stims-16509961009432955.csv"
$wsTOL.Range("B7").Value = "ZM_stims-1650996100927329.csv"

# Sheet5: vSAT_TO - update B2:B5
$wsvSAT.Range("B2").Value = "vSAT_stims-1650996100999307.csv"
$wsvSAT.Range("B3").Value = "vSAT_stims-16509961009833267.csv"
$wsvSAT.Range("B4").Value = "SAT_stims-16509961009432955.csv"
$wsvSAT.Range("B5").Value = "SAT_stims-16509961009673352.csv"
